$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 23

# Column A holds the date as plain text (matching the existing rows), so we
# temporarily force a text format before assigning the value to stop Excel
# from auto-converting the "MM/DD/YYYY" looking string into a date serial,
# then restore the cell to the default (unstyled) look, just like the rest
# of the date column.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "12/17/2025"
$ws.Cells.Item($row, 1).Style = "Normal"
$ws.Cells.Item($row, 2).Value = 11777.44
$ws.Cells.Item($row, 3).Value = 0.20854811384464
$ws.Cells.Item($row, 4).Value = 0.79145188615536
$ws.Cells.Item($row, 5).Value = -147.47
$ws.Cells.Item($row, 6).Value = -32.06
$ws.Cells.Item($row, 7).Value = -21195.48
$ws.Cells.Item($row, 8).Value = -69.45999999999999
$ws.Cells.Item($row, 9).Value = -496.69
$ws.Cells.Item($row, 10).Value = -16.82
